# Fri, May 22, 2020  3:05:36 PM
#
# Re-apply the table style used by the three data tables in this deck:
# switch them from the custom table style
# {8EDEF95A-6CD9-437D-A210-A1221DA6AC18} to the (built-in) style
# {52F8B93B-D0F2-4F8F-B91D-B8C6DC318CBD}.
#
# PowerPoint's object model does not let a table's style be assigned
# through the Style property directly (it raises "Table styles cannot
# be assigned through a property"); Table.ApplyStyle(styleId) is the
# supported call.

$oldStyleId = "{8EDEF95A-6CD9-437D-A210-A1221DA6AC18}"
$newStyleId = "{52F8B93B-D0F2-4F8F-B91D-B8C6DC318CBD}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
